# Generate Report for Handback
# -----------------------------------------------------------------------
# This script updates the localization-status workbook to reflect that
# the handback xliff/markdown files have now been generated:
#   * The overall status text moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is used.
#   * The "Latest Target File" / "Latest Handback File" / "Latest
#     Handback DateTime" columns on the per-locale sheets are filled in
#     (they were placeholders before).
#   * A couple of columns are widened so the newly-populated long file
#     names remain readable.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared by the Overview sheet (zh-cn/de-de status
#    columns) and by the Status column on each per-locale sheet.
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) zh-cn sheet: populate Latest Target File (I), Latest Handback File
#    (J) and Latest Handback DateTime (K) for both rows, and hyperlink
#    the target-file cells the same way column A is hyperlinked.
# ---------------------------------------------------------------------
$md1 = "082f1594-bb1a-4914-aa45-d4c54e0c230c.md"
$md2 = "e4fe09eb-9d83-4633-ade2-c7a1f08e1489.md"
$url1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78296da06db1b9f7e2d2672a894b70fe5628f0c5/e2e/082f1594-bb1a-4914-aa45-d4c54e0c230c.md"
$url2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/78296da06db1b9f7e2d2672a894b70fe5628f0c5/e2e/e4fe09eb-9d83-4633-ade2-c7a1f08e1489.md"

# Re-create every hyperlink on this sheet so that the new ones land in
# the same ref order (A2, I2, A3, I3) as the handoff links.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $url1, "", "", $md1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $url1, "", "", $md1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $url2, "", "", $md2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $url2, "", "", $md2)

$wsZhCn.Range("J2").Value = "082f1594-bb1a-4914-aa45-d4c54e0c230c.3510494cf7d55678270257f1e7552604e4d74714.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "e4fe09eb-9d83-4633-ade2-c7a1f08e1489.f462688ada0e6ec8618a91f7a6874f40b9181f51.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-29 21:01:55"
$wsZhCn.Range("K3").Value = "2016-08-29 21:01:55"

# ---------------------------------------------------------------------
# 3) de-de sheet: same treatment, with its own xliff file names and a
#    later handback timestamp.
# ---------------------------------------------------------------------
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $url1, "", "", $md1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $url1, "", "", $md1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $url2, "", "", $md2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $url2, "", "", $md2)

$wsDeDe.Range("J2").Value = "082f1594-bb1a-4914-aa45-d4c54e0c230c.3510494cf7d55678270257f1e7552604e4d74714.de-de.xlf"
$wsDeDe.Range("J3").Value = "e4fe09eb-9d83-4633-ade2-c7a1f08e1489.f462688ada0e6ec8618a91f7a6874f40b9181f51.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-29 21:02:12"
$wsDeDe.Range("K3").Value = "2016-08-29 21:02:12"

# ---------------------------------------------------------------------
# 4) Column widths: widen the columns that now hold long status text /
#    file names so the handback information stays legible.
#    (ColumnWidth is Excel's pixel-quantized "characters" unit; the
#    values below are the closest attainable approximation of the
#    ~29.98-character target width.)
# ---------------------------------------------------------------------
$wideWidth = 29.09
$fullWidth = 39.09

$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth   # E
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth   # F

$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth       # C (Status)
$wsZhCn.Columns.Item(9).ColumnWidth = $fullWidth        # I (Latest Target File)
$wsZhCn.Columns.Item(10).ColumnWidth = $fullWidth       # J (Latest Handback File)

$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth       # C (Status)
$wsDeDe.Columns.Item(9).ColumnWidth = $fullWidth        # I (Latest Target File)
$wsDeDe.Columns.Item(10).ColumnWidth = $fullWidth       # J (Latest Handback File)
